$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-23 Thursday" "2023-11-24 Friday"

Replace-Text "89÷9=9, 8" "58÷3=19, 1"
Replace-Text "32÷3=10, 2" "96÷3=32, 0"
Replace-Text "37÷7=5, 2" "92÷3=30, 2"
Replace-Text "45÷4=11, 1" "58÷4=14, 2"
Replace-Text "53÷6=8, 5" "21÷5=4, 1"

Replace-Text "76÷4=19, 0" "10÷2=5, 0"
Replace-Text "80÷5=16, 0" "71÷2=35, 1"
Replace-Text "75÷2=37, 1" "64÷2=32, 0"
Replace-Text "52÷9=5, 7" "39÷2=19, 1"
Replace-Text "44÷4=11, 0" "73÷4=18, 1"

Replace-Text "96÷2=48, 0" "80÷3=26, 2"
Replace-Text "92÷4=23, 0" "57÷2=28, 1"
Replace-Text "71÷6=11, 5" "37÷6=6, 1"
Replace-Text "38÷7=5, 3" "61÷3=20, 1"
Replace-Text "39÷3=13, 0" "93÷5=18, 3"

Replace-Text "43÷7=6, 1" "75÷3=25, 0"
Replace-Text "71÷5=14, 1" "98÷5=19, 3"
Replace-Text "13÷9=1, 4" "67÷9=7, 4"
Replace-Text "76÷5=15, 1" "18÷7=2, 4"
Replace-Text "30÷4=7, 2" "77÷5=15, 2"

Replace-Text "94÷8=11, 6" "50÷6=8, 2"
Replace-Text "21÷9=2, 3" "49÷6=8, 1"
Replace-Text "57÷4=14, 1" "51÷2=25, 1"
Replace-Text "94÷5=18, 4" "40÷8=5, 0"
Replace-Text "92÷5=18, 2" "13÷3=4, 1"
